$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.092.14"
$ws.Range("E2").Value = "  +4.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.057.45"
$ws.Range("E3").Value = "  +7.47%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.68"
$ws.Range("E5").Value = "  +8.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.51"
$ws.Range("E6").Value = "  +8.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.433"
$ws.Range("E8").Value = "  +6.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.20"
$ws.Range("E9").Value = "  +4.59%  "

$ws.Range("E10").Value = "  +8.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.371"
$ws.Range("E11").Value = "  +10.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.553.36"
$ws.Range("E12").Value = "  +5.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.127"
$ws.Range("E13").Value = "  +3.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.36"
$ws.Range("E14").Value = "  +2.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000164"
$ws.Range("E15").Value = "  +9.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.085.63"
$ws.Range("E16").Value = "  +4.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.038.90"
$ws.Range("E17").Value = "  +6.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.97"
$ws.Range("E18").Value = "  +3.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.26"
$ws.Range("E19").Value = "  +11.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.16"
$ws.Range("E20").Value = "  +10.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "336.88"
$ws.Range("E21").Value = "  +11.64%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.507"
$ws.Range("E23").Value = "  +9.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.98"
$ws.Range("E24").Value = "  +8.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.154.77"
$ws.Range("E25").Value = "  +5.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +7.82%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0935"
$ws.Range("E28").Value = "  +19.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.47"
$ws.Range("E29").Value = "  +5.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  +3.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  +8.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.81"
$ws.Range("E32").Value = "  +8.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  +8.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "154.02"
$ws.Range("E34").Value = "  +5.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.57"
$ws.Range("E35").Value = "  +8.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.42"
$ws.Range("E36").Value = "  +21.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.88"
$ws.Range("E37").Value = "  +9.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.24"
$ws.Range("E38").Value = "  +8.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0676"
$ws.Range("E39").Value = "  +7.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.081.05"
$ws.Range("E40").Value = "  +6.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.01"
$ws.Range("E41").Value = "  +5.16%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("E42").Value = "  +10.56%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.663"
$ws.Range("E44").Value = "  +9.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.228.12"
$ws.Range("E45").Value = "  +8.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0250"
$ws.Range("E46").Value = "  +14.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.37"
$ws.Range("E47").Value = "  +7.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.945"
$ws.Range("E48").Value = "  +10.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.03"
$ws.Range("E49").Value = "  +11.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.90"
$ws.Range("E50").Value = "  +4.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0864"
$ws.Range("E51").Value = "  +6.44%  "

